# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# Column G ("K") values were recomputed from the underlying per-game pitching
# log (replacing the previous "Strike#" derived figures) and rewritten back
# into the worksheet. The recalculated values for each row are applied here.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 2
    3  = 4
    4  = 1
    5  = 0
    6  = 0
    7  = 0
    8  = 2
    9  = 1
    10 = 0
    11 = 1
    12 = 2
    13 = 1
    14 = 2
    15 = 0
    16 = 2
    17 = 0
    18 = 2
    19 = 0
    20 = 2
    21 = 0
    22 = 1
    23 = 3
    24 = 0
    25 = 3
    26 = 2
    27 = 0
    28 = 1
    29 = 1
    30 = 1
    31 = 2
    32 = 1
    33 = 1
    34 = 1
    35 = 2
    36 = 2
    37 = 1
    38 = 1
    39 = 1
    40 = 1
    41 = 2
    42 = 1
    43 = 2
    44 = 2
    45 = 1
    46 = 0
    47 = 2
    48 = 2
    49 = 2
    50 = 2
    51 = 1
    52 = 1
    53 = 1
    54 = 3
    55 = 2
    56 = 2
    57 = 0
    58 = 1
    59 = 0
    60 = 1
    61 = 1
    62 = 1
    63 = 1
    64 = 2
    65 = 1
    66 = 1
    67 = 0
    68 = 1
    69 = 1
    70 = 0
    71 = 3
    72 = 1
    73 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
